$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 517.3333
$ws.Range("I29").Value = 26
$ws.Range("K29").Value = 78
$ws.Range("M29").Value = 203
$ws.Range("H32").Value = 853.36365
$ws.Range("I32").Value = 798.75
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 798.75
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -472.75
$ws.Range("N32").Value = -1651
$ws.Range("H41").Value = 231.55556
$ws.Range("I41").Value = 185.83333
$ws.Range("K41").Value = 185.83333
$ws.Range("M41").Value = 254.16667
$ws.Range("H55").Value = 443.66666
$ws.Range("I55").Value = 40
$ws.Range("J55").Value = 524.4
$ws.Range("K55").Value = 40
$ws.Range("L55").Value = 524.4
$ws.Range("M55").Value = 174
$ws.Range("N55").Value = -952.4
$ws.Range("H58").Value = 1705.909
$ws.Range("I58").Value = 1212.25
$ws.Range("J58").Value = 1988
$ws.Range("K58").Value = 3636.75
$ws.Range("L58").Value = 5964
$ws.Range("M58").Value = -3486.75
$ws.Range("N58").Value = -6264
$ws.Range("H62").Value = 3853.2856
$ws.Range("I62").Value = 3540.6365
$ws.Range("K62").Value = 3540.6365
$ws.Range("M62").Value = -2916.6365
$ws.Range("H65").Value = 3853.2856
$ws.Range("I65").Value = 3540.6365
$ws.Range("K65").Value = 17703.1825
$ws.Range("M65").Value = -14583.1825
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H129").Value = 1635.75
$ws.Range("I129").Value = 817.6
$ws.Range("J129").Value = 2999.3333
$ws.Range("K129").Value = 2452.8
$ws.Range("L129").Value = 8997.999899999999
$ws.Range("M129").Value = 2547.2
$ws.Range("N129").Value = -18997.9999
$ws.Range("H132").Value = 2257.889
$ws.Range("I132").Value = 2257.889
$ws.Range("K132").Value = 6773.667
$ws.Range("M132").Value = -4243.667
$ws.Range("H137").Value = 1890.6
$ws.Range("J137").Value = 2668.5
$ws.Range("L137").Value = 8005.5
$ws.Range("N137").Value = -13105.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 957
$ws.Range("I2").Value = 1028.6
$ws.Range("J2").Value = 599
$ws.Range("K2").Value = 1028.6
$ws.Range("L2").Value = 599
$ws.Range("M2").Value = -915.5999999999999
$ws.Range("N2").Value = -825
$ws.Range("H32").Value = 1694.0422
$ws.Range("I32").Value = 1575.3857
$ws.Range("K32").Value = 1575.3857
$ws.Range("M32").Value = -1288.3857
$ws.Range("H76").Value = 39983.332
$ws.Range("J76").Value = 39983.332
$ws.Range("L76").Value = 39983.332
$ws.Range("N76").Value = -40659.332
$ws.Range("H79").Value = 39983.332
$ws.Range("J79").Value = 39983.332
$ws.Range("L79").Value = 39983.332
$ws.Range("N79").Value = -42323.332
$ws.Range("H92").Value = 32516.334
$ws.Range("J92").Value = 32516.334
$ws.Range("L92").Value = 32516.334
$ws.Range("N92").Value = -37508.334
$ws.Range("H102").Value = 2073.7778
$ws.Range("J102").Value = 2694.25
$ws.Range("L102").Value = 2694.25
$ws.Range("N102").Value = -5938.25
$ws.Range("H116").Value = 957
$ws.Range("I116").Value = 1028.6
$ws.Range("J116").Value = 599
$ws.Range("K116").Value = 1028.6
$ws.Range("L116").Value = 599
$ws.Range("M116").Value = 1265.4
$ws.Range("N116").Value = -5187

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 957
$ws.Range("I3").Value = 1028.6
$ws.Range("J3").Value = 599
$ws.Range("K3").Value = 1028.6
$ws.Range("L3").Value = 599
$ws.Range("M3").Value = -914.5999999999999
$ws.Range("N3").Value = -827
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H94").Value = 1007.95
$ws.Range("I94").Value = 1048.1765
$ws.Range("K94").Value = 1048.1765
$ws.Range("M94").Value = -597.1765
$ws.Range("H100").Value = 13631.75
$ws.Range("J100").Value = 13631.75
$ws.Range("L100").Value = 13631.75
$ws.Range("N100").Value = -15795.75
$ws.Range("H105").Value = 3856
$ws.Range("I105").Value = 3692.7144
$ws.Range("K105").Value = 3692.7144
$ws.Range("M105").Value = -1945.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2399.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2399.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H86").Value = 7666.4287
$ws.Range("J86").Value = 8340.799999999999
$ws.Range("L86").Value = 8340.799999999999
$ws.Range("N86").Value = -10586.8
$ws.Range("H89").Value = 7666.4287
$ws.Range("J89").Value = 8340.799999999999
$ws.Range("L89").Value = 41704
$ws.Range("N89").Value = -52936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7286.75
$ws.Range("J80").Value = 6415.8335
$ws.Range("L80").Value = 19247.5005
$ws.Range("N80").Value = -21119.5005
$ws.Range("H83").Value = 7286.75
$ws.Range("J83").Value = 6415.8335
$ws.Range("L83").Value = 57742.5015
$ws.Range("N83").Value = -67102.5015
$ws.Range("H107").Value = 736.36365
$ws.Range("I107").Value = 680.6
$ws.Range("J107").Value = 782.8333
$ws.Range("K107").Value = 2041.8
$ws.Range("L107").Value = 2348.4999
$ws.Range("M107").Value = -121.8000000000002
$ws.Range("N107").Value = -6188.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 38333.332
$ws.Range("J18").Value = 38333.332
$ws.Range("L18").Value = 38333.332
$ws.Range("N18").Value = -38919.332
$ws.Range("H21").Value = 12513000
$ws.Range("J21").Value = 26000
$ws.Range("L21").Value = 26000
$ws.Range("N21").Value = -26346
$ws.Range("H30").Value = 12513000
$ws.Range("J30").Value = 26000
$ws.Range("L30").Value = 26000
$ws.Range("N30").Value = -26210
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 3000
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 15000
$ws.Range("N83").Value = -24984
$ws.Range("H113").Value = 4999
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -2829
$ws.Range("H122").Value = 4732.3335
$ws.Range("I122").Value = 4597.5
$ws.Range("J122").Value = 5002
$ws.Range("K122").Value = 13792.5
$ws.Range("L122").Value = 15006
$ws.Range("M122").Value = -11342.5
$ws.Range("N122").Value = -19906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3341.7896
$ws.Range("I46").Value = 2220
$ws.Range("J46").Value = 3742.4285
$ws.Range("K46").Value = 2220
$ws.Range("L46").Value = 3742.4285
$ws.Range("M46").Value = -2032
$ws.Range("N46").Value = -4118.4285
$ws.Range("H55").Value = 207.26666
$ws.Range("I55").Value = 190.66667
$ws.Range("K55").Value = 190.66667
$ws.Range("M55").Value = -17.66667000000001
$ws.Range("H127").Value = 75998
$ws.Range("J127").Value = 75998
$ws.Range("L127").Value = 75998
$ws.Range("N127").Value = -85918
$ws.Range("H136").Value = 2970
$ws.Range("I136").Value = 2782.2727
$ws.Range("K136").Value = 8346.8181
$ws.Range("M136").Value = -5796.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 34750
$ws.Range("J45").Value = 44500
$ws.Range("L45").Value = 44500
$ws.Range("N45").Value = -45482
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H126").Value = 1140.5714
$ws.Range("I126").Value = 496.25
$ws.Range("K126").Value = 1488.75
$ws.Range("M126").Value = 981.25
